$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skills")

# --- New "Mage" section title row (row 11), copied formatting from the
# existing "Knight" title row (row 1) but with a new accent fill color ---
$ws.Range("A1:G1").Copy()
$ws.Range("A11:G11").PasteSpecial(-4122)
$ws.Range("A11:G11").Interior.Color = 13307391
$ws.Range("A11").Value = "Mage"
$ws.Range("A11:G11").Merge()

# --- Column header row (row 12), same layout/style as row 2 ---
$ws.Range("A2:G2").Copy()
$ws.Range("A12:G12").PasteSpecial(-4122)
$ws.Range("A12").Value = "ID"
$ws.Range("B12").Value = "Name"
$ws.Range("C12").Value = "Description"
$ws.Range("D12").Value = "Type"
$ws.Range("E12").Value = "Stamina"
$ws.Range("F12").Value = "Base damage"
$ws.Range("G12").Value = "Cooldown"

# --- Skill rows 13-19, same layout/style as rows 3-9 ---
$ws.Range("A3:G9").Copy()
$ws.Range("A13:G19").PasteSpecial(-4122)

$ws.Range("A13").Value = 7
$ws.Range("B13").Value = "Mage"
$ws.Range("C13").Value = "Casts fireball at a target dealing instant damage, plus setting the target on  fire dealing reoccurring damage for 5 seconds. Instant"
$ws.Range("D13").Value = "Damaging"
$ws.Range("E13").Value = 1000
$ws.Range("F13").Value = 100
$ws.Range("G13").Value = 2

$ws.Range("A14").Value = 8
$ws.Range("A15").Value = 9
$ws.Range("A16").Value = 10
$ws.Range("A17").Value = 11
$ws.Range("A18").Value = 12
$ws.Range("A19").Value = 13

# --- Column C needs to be wide enough to fit the new longer description ---
$ws.Columns.Item(3).ColumnWidth = 116.6

# --- Data validation dropdown (Type column) needs to cover the new rows ---
$ws.Range("D13:D19").Validation.Add(3, 1, 1, "=Variables!$A$2:$A$5")

# --- Print orientation ---
$ws.PageSetup.Orientation = 1

# --- Final selection, matches where the author ended up ---
$ws.Range("C26").Select()
